# Horarios actualizados Linea 141 - 875
# Updates the scrape snapshot across the three route sheets:
#   Sheet 1 "LP1912"     (full line)
#   Sheet 2 "LP1912-215" (215 sub-line)
#   Sheet 3 "6203-6173"  (previously empty feed, now has its first row)

$wb = $excel.ActiveWorkbook

$NEW_SCRAPE = "03:52:04"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $NEW_SCRAPE"
$ws1.Range("A3").Value = "Total filas: 7"

$sheet1Rows = @(
    @("04:01", "81_EL PELIGRO", 9,   "LP1912"),
    @("04:46", "215A_EL PATO",  54,  "LP1912"),
    @("04:53", "11_ETCHEVERRY", 61,  "LP1912"),
    @("05:16", "17_ROMERO",     84,  "LP1912"),
    @("05:22", "23_HERNANDEZ",  90,  "LP1912"),
    @("05:35", "215B_EL PATO",  103, "LP1912"),
    @("05:46", "15_ABASTO",     114, "LP1912")
)

$r = 6
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $NEW_SCRAPE
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $NEW_SCRAPE"
$ws2.Range("A3").Value = "Total filas: 2"

$sheet2Rows = @(
    @("04:46", "215A_EL PATO", 54,  "LP1912"),
    @("05:35", "215B_EL PATO", 103, "LP1912")
)

$r = 6
foreach ($row in $sheet2Rows) {
    $ws2.Cells.Item($r, 1).Value = $NEW_SCRAPE
    $ws2.Cells.Item($r, 2).Value = $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173 (was empty, now gets its first data row + headers)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $NEW_SCRAPE"
$ws3.Range("A3").Value = "Total filas: 1"

$ws3.Range("A5").Value = "Hora_Scrap"
$ws3.Range("B5").Value = "Hora_Llegada"
$ws3.Range("C5").Value = "Linea"
$ws3.Range("D5").Value = "Minutos"
$ws3.Range("E5").Value = "Parada"

$ws3.Range("A6").Value = $NEW_SCRAPE
$ws3.Range("B6").Value = "05:44"
$ws3.Range("C6").Value = "215A_LA PLATA"
$ws3.Range("D6").Value = 112
$ws3.Range("E6").Value = "L6173"
